$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.688085675239563
$ws.Range("B1").Value = 1.941734910011292
$ws.Range("C1").Value = 1.99854838848114
$ws.Range("D1").Value = 2.470536231994629
$ws.Range("E1").Value = 3.371030569076538
